$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16 (for the D3,D4,D5 diode entry), shifting everything
# below down by one. This also pushes the PS1/PS2/PS3 power-supply rows from
# 31-33 down to 32-34.
$ws.Rows("16:16").Insert()

# Fill in the new diode row (1N5819 Schottky diodes, replacing the old DC/DC
# converters used elsewhere in the BOM).
$ws.Range("A16").Value = "D3, D4,D5"
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = "1N5819"
$ws.Range("D16").Value = "SOD-123"
$ws.Range("E16").Value = "621-1N5819HW-F "
$ws.Range("E16").WrapText = $true

# Replace the old PXO7806/PXO7803/PXO7805 DC/DC converters (PS1, PS2, PS3)
# with the new R-78Kx regulators (U17, U18, U19).
$ws.Range("A32").Value = "U17"
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = "R-78K6.5-0.5 "
$ws.Range("E32").Value = "919-R-78K6.5-0.5 "

$ws.Range("A33").Value = "U18"
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = "R-78K5.0-0.5 "
$ws.Range("E33").Value = "919-R-78K5.0-0.5 "

$ws.Range("A34").Value = "U19"
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = "R-78K63.3-0.5 "
$ws.Range("E34").Value = "919-R-78K3.3-0.5 "

# Match the author's final cursor position recorded in the saved file.
$ws.Range("D16").Select()
